# task-list.xlsx update
# Adds an "Owner" column value to several rows (Mobile sheet col D, Web App
# sheet col E), renames two Mobile-sheet task names, and appends a few new
# task rows to both the "Mobile" and "Web App" sheets.

$wb = $excel.ActiveWorkbook

$wsMobile = $wb.Worksheets.Item("Mobile")
$wsWebApp = $wb.Worksheets.Item("Web App")

# ---------------------------------------------------------------------
# Mobile sheet ("sheet1")
# ---------------------------------------------------------------------

# Rename two existing tasks
$wsMobile.Cells.Item(9, 2).Value  = "timeline page"
$wsMobile.Cells.Item(10, 2).Value = "Sync page"

# Fill in the new "Owner" column (D) for existing rows
$wsMobile.Cells.Item(7, 4).Value  = "Megha"
$wsMobile.Cells.Item(8, 4).Value  = "Rajat"
$wsMobile.Cells.Item(9, 4).Value  = "Eric"
$wsMobile.Cells.Item(10, 4).Value = "Pranav"
$wsMobile.Cells.Item(11, 4).Value = "Eric"

# Append 3 new task rows
$wsMobile.Cells.Item(13, 1).Value = 12
$wsMobile.Cells.Item(13, 2).Value = "Jobs Page"
$wsMobile.Cells.Item(13, 3).Value = "On the day"
$wsMobile.Cells.Item(13, 4).Value = "Eric"

$wsMobile.Cells.Item(14, 1).Value = 13
$wsMobile.Cells.Item(14, 2).Value = "Jobs Details"
$wsMobile.Cells.Item(14, 3).Value = "On the day"
$wsMobile.Cells.Item(14, 4).Value = "Rajat"

$wsMobile.Cells.Item(15, 1).Value = 14
$wsMobile.Cells.Item(15, 2).Value = "My Referals"
$wsMobile.Cells.Item(15, 3).Value = "On the day"
$wsMobile.Cells.Item(15, 4).Value = "Megha"

# ---------------------------------------------------------------------
# Web App sheet ("sheet2")
# ---------------------------------------------------------------------

# Fix/fill the "Owner" column (E) for existing rows
$wsWebApp.Cells.Item(7, 5).Value  = "Richa"
$wsWebApp.Cells.Item(8, 5).Value  = "Sadhvi"
$wsWebApp.Cells.Item(9, 5).Value  = "Sadhvi"
$wsWebApp.Cells.Item(10, 5).Value = "Richa"
$wsWebApp.Cells.Item(13, 5).Value = "keshav"
$wsWebApp.Cells.Item(14, 5).Value = "keshav"
$wsWebApp.Cells.Item(15, 5).Value = "Nitin"
$wsWebApp.Cells.Item(16, 5).Value = "Nitin"
$wsWebApp.Cells.Item(18, 5).Value = "Keshav"
$wsWebApp.Cells.Item(28, 5).Value = "Keshav"
$wsWebApp.Cells.Item(29, 5).Value = "Nitin"

# Append 1 new task row
$wsWebApp.Cells.Item(30, 1).Value = 29
$wsWebApp.Cells.Item(30, 2).Value = "Login Page"

# ---------------------------------------------------------------------
# Selection / active-cell bookkeeping (matches the saved view state)
# ---------------------------------------------------------------------

$wsWebApp.Activate()
$wsWebApp.Range("E24").Select()
$excel.ActiveWindow.ScrollRow = 11

$wsDb = $wb.Worksheets.Item("Database Design")
$wsDb.Activate()
$wsDb.Range("D6").Select()

$wsMobile.Activate()
$wsMobile.Range("D5").Select()
